$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "runs" (C) and "balls" (D) values between row 3 and row 4.
# The source cells are stored as text, so prefix with an apostrophe to
# keep Excel from reinterpreting the values as numbers.
$ws.Range("C3").Value = "'2"
$ws.Range("D3").Value = "'2"
$ws.Range("C4").Value = "'5"
$ws.Range("D4").Value = "'5"
